$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @()
$rows += ,@(1169, 45534.5, 7.752, 7.86, 7.644, 7.759, 181530.37)
$rows += ,@(1170, 45534.58333333334, 7.767, 7.824, 7.383, 7.415, 349549.84)
$rows += ,@(1171, 45534.66666666666, 7.419, 7.489, 7.336, 7.439, 221138.81)
$rows += ,@(1172, 45534.75, 7.441, 7.62, 7.423, 7.537, 131897.86)
$rows += ,@(1173, 45534.83333333334, 7.533, 7.599, 7.52, 7.582, 42539.73)
$rows += ,@(1174, 45534.91666666666, 7.585, 7.656, 7.581, 7.617, 70244.44)
$rows += ,@(1175, 45535, 7.618, 7.668, 7.595, 7.627, 100875.78)
$rows += ,@(1176, 45535.08333333334, 7.626, 7.679, 7.622, 7.654, 73785.12)
$rows += ,@(1177, 45535.16666666666, 7.654, 7.722, 7.654, 7.66, 63851.77)
$rows += ,@(1178, 45535.25, 7.66, 7.699, 7.611, 7.618, 51698.88)
$rows += ,@(1179, 45535.33333333334, 7.618, 7.631, 7.519, 7.523, 87657.78)
$rows += ,@(1180, 45535.41666666666, 7.512, 7.58, 7.503, 7.537, 85213.53)
$rows += ,@(1181, 45535.5, 7.54, 7.558, 7.461, 7.462, 50753.53)
$rows += ,@(1182, 45535.58333333334, 7.46, 7.539, 7.457, 7.537, 49009.9)
$rows += ,@(1183, 45535.66666666666, 7.537, 7.538, 7.435, 7.493, 41519.26)
$rows += ,@(1184, 45535.75, 7.493, 7.562, 7.478, 7.5, 47406.04)
$rows += ,@(1185, 45535.83333333334, 7.501, 7.53, 7.487, 7.523, 17235.76)
$rows += ,@(1186, 45535.91666666666, 7.55, 7.576, 7.536, 7.548, 7388.96)
$rows += ,@(1187, 45536, 7.56, 7.583, 7.453, 7.454, 40674.37)
$rows += ,@(1188, 45536.08333333334, 7.455, 7.487, 7.41, 7.481, 44923.04)
$rows += ,@(1189, 45536.16666666666, 7.482, 7.497, 7.351, 7.367, 61923.62)
$rows += ,@(1190, 45536.25, 7.365, 7.441, 7.351, 7.436, 57890.3)
$rows += ,@(1191, 45536.33333333334, 7.437, 7.449, 7.352, 7.354, 41074.14)
$rows += ,@(1192, 45536.41666666666, 7.357, 7.363, 7.286, 7.357, 124859.17)
$rows += ,@(1193, 45536.5, 7.359, 7.373, 7.258, 7.314, 65197.76)
$rows += ,@(1194, 45536.58333333334, 7.315, 7.436, 7.189, 7.402, 98533.97)
$rows += ,@(1195, 45536.66666666666, 7.404, 7.435, 7.245, 7.297, 82629.75999999999)
$rows += ,@(1196, 45536.75, 7.296, 7.451, 7.272, 7.43, 103607.69)
$rows += ,@(1197, 45536.83333333334, 7.425, 7.462, 7.32, 7.352, 28663.61)
$rows += ,@(1198, 45536.91666666666, 7.348, 7.361, 7.087, 7.191, 178476.37)
$rows += ,@(1199, 45537, 7.184, 7.241, 7.131, 7.189, 104285.05)
$rows += ,@(1200, 45537.08333333334, 7.192, 7.24, 7.162, 7.224, 54002.21)
$rows += ,@(1201, 45537.16666666666, 7.225, 7.311, 7.176, 7.276, 56550.87)
$rows += ,@(1202, 45537.25, 7.275, 7.302, 7.115, 7.146, 100812.48)
$rows += ,@(1203, 45537.33333333334, 7.145, 7.343, 7.144, 7.331, 176690.48)
$rows += ,@(1204, 45537.41666666666, 7.327, 7.427, 7.321, 7.36, 124424.3)
$rows += ,@(1205, 45537.5, 7.363, 7.38, 7.309, 7.32, 23494.28)
$rows += ,@(1206, 45537.58333333334, 7.293, 7.429, 7.238, 7.391, 140178.77)
$rows += ,@(1207, 45537.66666666666, 7.391, 7.435, 7.342, 7.399, 94488.92)
$rows += ,@(1208, 45537.75, 7.397, 7.427, 7.363, 7.387, 46983.78)
$rows += ,@(1209, 45537.83333333334, 7.393, 7.511, 7.351, 7.445, 107823.35)
$rows += ,@(1210, 45537.91666666666, 7.446, 7.552, 7.442, 7.521, 82835.38)
$rows += ,@(1211, 45538, 7.524, 7.568, 7.493, 7.567, 95273.28999999999)
$rows += ,@(1212, 45538.08333333334, 7.567, 7.645, 7.521, 7.534, 131912.49)
$rows += ,@(1213, 45538.16666666666, 7.533, 7.555, 7.47, 7.537, 46209.43)
$rows += ,@(1214, 45538.25, 7.539, 7.54, 7.469, 7.518, 76489.03999999999)
$rows += ,@(1215, 45538.33333333334, 7.517, 7.543, 7.421, 7.43, 111298.34)
$rows += ,@(1216, 45538.41666666666, 7.431, 7.488, 7.407, 7.475, 59323.11)
$rows += ,@(1217, 45538.5, 7.47, 7.519, 7.282, 7.313, 133348.14)
$rows += ,@(1218, 45538.58333333334, 7.314, 7.316, 7.149, 7.164, 198016.76)
$rows += ,@(1219, 45538.66666666666, 7.161, 7.259, 7.151, 7.22, 96049.97)
$rows += ,@(1220, 45538.75, 7.22, 7.269, 7.173, 7.252, 39097.89)
$rows += ,@(1221, 45538.83333333334, 7.249, 7.293, 7.231, 7.257, 25104.14)
$rows += ,@(1222, 45538.91666666666, 7.258, 7.269, 7.144, 7.162, 44306.53)
$rows += ,@(1223, 45539, 7.158, 7.22, 6.9, 7.103, 331155.3)
$rows += ,@(1224, 45539.08333333334, 7.102, 7.196, 7.086, 7.191, 104903.88)
$rows += ,@(1225, 45539.16666666666, 7.191, 7.23, 7.16, 7.179, 106549.64)
$rows += ,@(1226, 45539.25, 7.178, 7.279, 7.166, 7.253, 128307.34)
$rows += ,@(1227, 45539.33333333334, 7.255, 7.31, 7.202, 7.234, 104384.55)
$rows += ,@(1228, 45539.41666666666, 7.236, 7.255, 7.151, 7.217, 128333.2)
$rows += ,@(1229, 45539.5, 7.216, 7.299, 7.171, 7.295, 147307.39)
$rows += ,@(1230, 45539.58333333334, 7.294, 7.573, 7.234, 7.562, 344462.55)
$rows += ,@(1231, 45539.66666666666, 7.562, 7.76, 7.539, 7.62, 237266.55)
$rows += ,@(1232, 45539.75, 7.621, 7.644, 7.522, 7.538, 83371.85000000001)
$rows += ,@(1233, 45539.83333333334, 7.555, 7.657, 7.549, 7.647, 67344.97)
$rows += ,@(1234, 45539.91666666666, 7.647, 7.687, 7.554, 7.557, 96646.19)
$rows += ,@(1235, 45540, 7.557, 7.645, 7.549, 7.626, 81664.10000000001)
$rows += ,@(1236, 45540.08333333334, 7.624, 7.64, 7.441, 7.466, 125382.61)
$rows += ,@(1237, 45540.16666666666, 7.466, 7.499, 7.387, 7.428, 283692.96)
$rows += ,@(1238, 45540.25, 7.426, 7.468, 7.362, 7.436, 110105.28)
$rows += ,@(1239, 45540.33333333334, 7.435, 7.464, 7.387, 7.41, 83340.89999999999)
$rows += ,@(1240, 45540.41666666666, 7.409, 7.437, 7.328, 7.333, 293616.86)
$rows += ,@(1241, 45540.5, 7.333, 7.335, 7.33, 7.333, 884.47)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $row = $rows[$i]
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}

# Copy the date/time style from the last original row (A1168) to the new column-A cells
$ws.Range("A1168").Copy()
$ws.Range("A1169:A1241").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Host "Inserted rows 1169-1241"
